$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c2 = $ws.Range("C2")
$c2.Borders.Item(10).LineStyle = -4142
$c2.Borders.Item(10).LineStyle = 1
$c2.Borders.Item(10).ThemeColor = 1
Write-Host "done"
